# Update Effort Estimation sheet and Cost Summary sheet data.
# The module/feature breakdown changed from an "Auth/Dashboard/Backend"
# style breakdown to a "User Authentication / Product Catalog" breakdown,
# the table shrank from 14 data rows to 11, and the derived pricing
# numbers in the Cost Summary sheet were recalculated to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Effort Estimation")
$ws2 = $wb.Worksheets.Item("Cost Summary")

# ---------------------------------------------------------------------
# Sheet 1: "Effort Estimation"
# ---------------------------------------------------------------------

# First, clear out the old rows 2-14 entirely so no stale cells remain
# from the previous (longer) table before writing the new 2-11 block.
$ws1.Range("A2:I14").ClearContents()

# Row 2: User Authentication / User Registration / Frontend Implementation
$ws1.Cells.Item(2,1).Value = "User Authentication"
$ws1.Cells.Item(2,2).Value = "User Registration"
$ws1.Cells.Item(2,3).Value = "Frontend Implementation"
$ws1.Cells.Item(2,4).Value = 5
$ws1.Cells.Item(2,5).Value = 1
$ws1.Cells.Item(2,6).Value = 0.8999999999999999
$ws1.Cells.Item(2,7).Value = 7
$ws1.Cells.Item(2,8).Value = 1.4
$ws1.Cells.Item(2,9).Value = 1.26

# Row 3: User Authentication / User Registration / Email Verification
$ws1.Cells.Item(3,1).Value = "User Authentication"
$ws1.Cells.Item(3,2).Value = "User Registration"
$ws1.Cells.Item(3,3).Value = "Email Verification"
$ws1.Cells.Item(3,4).Value = 3
$ws1.Cells.Item(3,5).Value = 0.6000000000000001
$ws1.Cells.Item(3,6).Value = 0.54
$ws1.Cells.Item(3,7).Value = 4
$ws1.Cells.Item(3,8).Value = 0.8
$ws1.Cells.Item(3,9).Value = 0.72

# Row 4: User Authentication / Login / Frontend Implementation
$ws1.Cells.Item(4,1).Value = "User Authentication"
$ws1.Cells.Item(4,2).Value = "Login"
$ws1.Cells.Item(4,3).Value = "Frontend Implementation"
$ws1.Cells.Item(4,4).Value = 3
$ws1.Cells.Item(4,5).Value = 0.6000000000000001
$ws1.Cells.Item(4,6).Value = 0.54
$ws1.Cells.Item(4,7).Value = 4
$ws1.Cells.Item(4,8).Value = 0.8
$ws1.Cells.Item(4,9).Value = 0.72

# Row 5: User Authentication / Login / Session Management
$ws1.Cells.Item(5,1).Value = "User Authentication"
$ws1.Cells.Item(5,2).Value = "Login"
$ws1.Cells.Item(5,3).Value = "Session Management"
$ws1.Cells.Item(5,4).Value = 2
$ws1.Cells.Item(5,5).Value = 0.4
$ws1.Cells.Item(5,6).Value = 0.36
$ws1.Cells.Item(5,7).Value = 3
$ws1.Cells.Item(5,8).Value = 0.6000000000000001
$ws1.Cells.Item(5,9).Value = 0.54

# Row 6: Product Catalog / Product Listing / Frontend Implementation
$ws1.Cells.Item(6,1).Value = "Product Catalog"
$ws1.Cells.Item(6,2).Value = "Product Listing"
$ws1.Cells.Item(6,3).Value = "Frontend Implementation"
$ws1.Cells.Item(6,4).Value = 7
$ws1.Cells.Item(6,5).Value = 1.4
$ws1.Cells.Item(6,6).Value = 1.26
$ws1.Cells.Item(6,7).Value = 10
$ws1.Cells.Item(6,8).Value = 2
$ws1.Cells.Item(6,9).Value = 1.8

# Row 7: Product Catalog / Product Listing / Filtering and Sorting
$ws1.Cells.Item(7,1).Value = "Product Catalog"
$ws1.Cells.Item(7,2).Value = "Product Listing"
$ws1.Cells.Item(7,3).Value = "Filtering and Sorting"
$ws1.Cells.Item(7,4).Value = 5
$ws1.Cells.Item(7,5).Value = 1
$ws1.Cells.Item(7,6).Value = 0.8999999999999999
$ws1.Cells.Item(7,7).Value = 7
$ws1.Cells.Item(7,8).Value = 1.4
$ws1.Cells.Item(7,9).Value = 1.26

# Row 8: Product Catalog / Product Details / Frontend Implementation
$ws1.Cells.Item(8,1).Value = "Product Catalog"
$ws1.Cells.Item(8,2).Value = "Product Details"
$ws1.Cells.Item(8,3).Value = "Frontend Implementation"
$ws1.Cells.Item(8,4).Value = 5
$ws1.Cells.Item(8,5).Value = 1
$ws1.Cells.Item(8,6).Value = 0.8999999999999999
$ws1.Cells.Item(8,7).Value = 7
$ws1.Cells.Item(8,8).Value = 1.4
$ws1.Cells.Item(8,9).Value = 1.26

# Row 9: Product Catalog / Product Details / Add to Cart
$ws1.Cells.Item(9,1).Value = "Product Catalog"
$ws1.Cells.Item(9,2).Value = "Product Details"
$ws1.Cells.Item(9,3).Value = "Add to Cart"
$ws1.Cells.Item(9,4).Value = 3
$ws1.Cells.Item(9,5).Value = 0.6000000000000001
$ws1.Cells.Item(9,6).Value = 0.54
$ws1.Cells.Item(9,7).Value = 4
$ws1.Cells.Item(9,8).Value = 0.8
$ws1.Cells.Item(9,9).Value = 0.72

# Row 10: Total row
$ws1.Cells.Item(10,1).Value = "Total"
$ws1.Cells.Item(10,3).Value = "Total"
$ws1.Cells.Item(10,4).Value = 33
$ws1.Cells.Item(10,5).Value = 6.6
$ws1.Cells.Item(10,6).Value = 5.94
$ws1.Cells.Item(10,7).Value = 46
$ws1.Cells.Item(10,8).Value = 9.200000000000001
$ws1.Cells.Item(10,9).Value = 8.279999999999999

# Row 11: Units row
$ws1.Cells.Item(11,3).Value = "Units"
$ws1.Cells.Item(11,4).Value = "days"
$ws1.Cells.Item(11,5).Value = "days"
$ws1.Cells.Item(11,6).Value = "days"
$ws1.Cells.Item(11,7).Value = "days"
$ws1.Cells.Item(11,8).Value = "days"
$ws1.Cells.Item(11,9).Value = "days"

# Remove the now-unused rows 12-14 entirely (the table shrank from
# A1:I14 to A1:I11).
$ws1.Range("A12:I14").Delete()

# ---------------------------------------------------------------------
# Sheet 2: "Cost Summary"
# ---------------------------------------------------------------------

# Column B narrows slightly (20.7109375 -> 19.7109375 character-widths).
$ws2.Columns.Item(2).ColumnWidth = 19.7109375

$ws2.Cells.Item(2,2).Value = 43.56
$ws2.Cells.Item(2,4).Value = "₹5,227.20"

$ws2.Cells.Item(3,2).Value = 60.72000000000001
$ws2.Cells.Item(3,4).Value = "₹7,772.16"

$ws2.Cells.Item(4,2).Value = 6.534000000000001
$ws2.Cells.Item(4,4).Value = "₹627.26"

$ws2.Cells.Item(5,4).Value = "₹13,626.62"
